$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.272.68"
$ws.Range("D3").Value = "3.495.49"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "4.090.94"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "3.493.85"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "64.346.76"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.571"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").Value = "3.635.02"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("E32").Value = "  -6.07%  "
$ws.Range("D33").Value = "3.517.09"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("E34").Value = "  +4.07%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0779"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.806"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("E47").Value = "  -4.51%  "
$ws.Range("D48").Value = "2.466.15"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.892"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -1.20%  "
